# Weekly update: insert two new report rows (row 41 and 42) for
# "Hortaliza, Terminal La Palmera de La Serena - Jengibre", pushing the
# previously existing rows 41-90 down to rows 43-92.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 41-90 down by two rows, opening up two blank
# rows at 41 and 42 (formats/number formats are carried along by Excel).
$ws.Rows("41:42").Insert(-4121)

# Row 41: new weekly entry
$ws.Cells.Item(41, 1).Value  = 8
$ws.Cells.Item(41, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(41, 3).Value  = "Coquimbo"
$ws.Cells.Item(41, 4).Value  = 44923
$ws.Cells.Item(41, 5).Value  = 4
$ws.Cells.Item(41, 6).Value  = 100114007
$ws.Cells.Item(41, 7).Value  = "Jengibre"
$ws.Cells.Item(41, 8).Value  = "Sin especificar"
$ws.Cells.Item(41, 9).Value  = "Primera"
$ws.Cells.Item(41, 10).Value = 450
$ws.Cells.Item(41, 11).Value = 14000
$ws.Cells.Item(41, 12).Value = 15000
$ws.Cells.Item(41, 13).Value = 14500
$ws.Cells.Item(41, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(41, 15).Value = "Perú"
$ws.Cells.Item(41, 16).Value = 1115
$ws.Cells.Item(41, 17).Value = 13
$ws.Cells.Item(41, 18).Value = "Hortaliza"

# Row 42: new weekly entry
$ws.Cells.Item(42, 1).Value  = 8
$ws.Cells.Item(42, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(42, 3).Value  = "Coquimbo"
$ws.Cells.Item(42, 4).Value  = 44923
$ws.Cells.Item(42, 5).Value  = 4
$ws.Cells.Item(42, 6).Value  = 100114007
$ws.Cells.Item(42, 7).Value  = "Jengibre"
$ws.Cells.Item(42, 8).Value  = "Sin especificar"
$ws.Cells.Item(42, 9).Value  = "Primera"
$ws.Cells.Item(42, 10).Value = 450
$ws.Cells.Item(42, 11).Value = 14000
$ws.Cells.Item(42, 12).Value = 15000
$ws.Cells.Item(42, 13).Value = 14500
$ws.Cells.Item(42, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(42, 15).Value = "Perú"
$ws.Cells.Item(42, 16).Value = 1115
$ws.Cells.Item(42, 17).Value = 13
$ws.Cells.Item(42, 18).Value = "Hortaliza"
